$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- New "Custom name" column (H) ---
# Header cell H2: copy the format of the neighbouring header (G2), keep that
# formatting, and give it its own text/shared-string entry.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$ws.Range("H2").Value = "Custom name"

# Data cells H3:H12: copy the format of the neighbouring data column (E3) so
# the new column matches the look of D:G; leave the cells themselves empty.
$ws.Range("E3").Copy()
$ws.Range("H3:H12").PasteSpecial($xlPasteFormats)

# --- Bugfix: clear the stray example values that used to live in D3:G3 ---
$ws.Range("D3:G3").ClearContents()

# --- Column widths / visibility ---
# Columns I (9) and J (10) are already hidden with width 0 in the source
# workbook, so they are intentionally left untouched here - touching them
# would just re-emit identical state while needlessly splitting the <col>
# range that currently spans columns 9:11 in the saved file.
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 26.33203125
$ws.Columns.Item(11).ColumnWidth = 10.5
$ws.Columns.Item(11).Hidden = $true

# --- Selection moves to the new column's first data cell ---
$ws.Range("H3").Select()
